# IST price update 2025-12-20 19:07
# Insert a new price-check column before column B (shifts old B:E -> C:F)
# and populate it with the 2025-12-21 00:32 timestamp header + prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the existing data column width so the freshly-inserted column matches it.
$dataColWidth = $ws.Range("C1").ColumnWidth

# Insert one new column at B; existing B:E (and their formatting) shift to C:F.
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("B1").ColumnWidth = $dataColWidth

# New timestamp header for the freshly inserted column.
$ws.Range("B1").Value = "2025-12-21 00:32"

# New price snapshot values for the inserted column.
$ws.Range("B2").Value = 929
$ws.Range("B3").Value = 569
$ws.Range("B4").Value = 299
$ws.Range("B5").Value = 569
$ws.Range("B6").Value = 499
$ws.Range("B7").Value = 569
$ws.Range("B8").Value = 929
$ws.Range("B9").Value = 299
$ws.Range("B10").Value = 299
$ws.Range("B11").Value = 929
$ws.Range("B12").Value = 569
$ws.Range("B13").Value = 569
$ws.Range("B14").Value = 499
$ws.Range("B15").Value = 499
$ws.Range("B16").Value = 299
$ws.Range("B17").Value = 929
$ws.Range("B18").Value = 499
$ws.Range("B19").Value = 1497
$ws.Range("B20").Value = 929
$ws.Range("B21").Value = 499
$ws.Range("B22").Value = 299
$ws.Range("B23").Value = 1299
$ws.Range("B24").Value = 929
$ws.Range("B25").Value = 929
$ws.Range("B26").Value = 1299
